$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy header style (B1, which has style index 1) into the new A1 cell ---
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# --- Copy header style into the new I1:K1 header cells ---
$ws.Range("B1").Copy()
$ws.Range("I1:K1").PasteSpecial(-4122)

# --- Row 1 header values ---
$ws.Range("A1").Value = "id"
$ws.Range("E1").Value = "pwy_ave_mcs"
$ws.Range("F1").Value = "mcs_rxn_1"
$ws.Range("G1").Value = "mcs_rxn_2"
$ws.Range("H1").Value = "mcs_rxn_3"
$ws.Range("I1").Value = "uniprot_rxn_1"
$ws.Range("J1").Value = "uniprot_rxn_2"
$ws.Range("K1").Value = "uniprot_rxn_3"

# --- Row 2 values ---
$ws.Range("A2").Value = 2
$ws.Range("D2").Value = 37.58993642927125
$ws.Range("E2").Value = 0.6044071669071669
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3839285714285714
$ws.Range("H2").Value = 0.7626262626262625
$ws.Range("I2").Value = "Q0C8L3, P54987, B3IUN8"
$ws.Range("J2").Value = "O87943, O87942, O87944"
$k2 = @"
Q5PDP3, Q989E7, Q3A133, Q9A5V1, P62456, Q02136, Q57MS3, Q606Q2, Q8XV79, Q46WL2, Q2P3K3, Q5LL27, Q82AA6, Q8EFB1, Q5WDH8, Q8ESR8, P16245, P62458, Q8YSM8, Q313T7, Q8YWL4, Q5QWP8, Q3Z878, P60862, Q4JW59, Q2RGV8, Q39K91, Q8P9P2, Q5ZW87, A0A0J9X7D2, Q97KI2, Q9K6Z2, P06988, Q47XB8, Q8PLG9, Q3MEV7, Q6D411, O34651, Q5V574, A1BPP9, Q722Y3, P60859, Q9PBC5, Q30TA4, Q8TXG3, P24226, Q7V5N9, Q92S26, Q9RQ88, Q5F7D8, Q8ZFX5, Q7VFF5, Q5HKN8, Q47AL5, P62457, Q39YP7, Q9JTH9, Q8ZY17, Q7NL02, O33775, Q6FEC8, Q323J2, Q5FS83, Q8R882, Q8F393, P62460, Q7M8K9, Q9HPW5, B0SS41, B0S9F2, Q65EF9, Q62GD9, Q5SKC1, Q3J7H1, Q4A044, Q63DX2, Q2YAU5, Q8D8Q0, Q318Q4, Q46J19, Q8G2R2, Q3AD53, Q6G5Z7, Q7WDY4, Q9P777, Q5E638, Q4FNE2, Q845V3, Q7VQX0, Q5KVC6, Q8G4S9, Q4QN74, Q8GKZ1, Q9ZHE6, Q3ICF0, Q3Z0G5, Q3J079, Q58851, O59626, Q48EC9, P44001, P63954, Q7W2Y4, P73058, Q9RSI4, Q87QL1, Q8X8T3, O26327, Q3ZXL7, Q81G06, Q8GDP4, P60861, Q4ZNV9, Q7V004, Q6HLE7, Q9RQ85, P18786, Q8DGR2, Q3B5E3, Q65RB1, P63953, Q8Z5K0, Q5X5W9, Q82W26, Q8UHX1, P9WNW8, Q5WX91, P60858, Q7VA26, Q8Y9G1, Q5HCL9, Q6GDC6, Q9CLM4, Q8ABA9, Q9C5U8, P59397, Q4J8I8, P10370, Q3SW26, Q5YYQ0, Q6A8L5, Q5MZ26, Q6AE76, Q8FG52, Q7U8K7, Q9L6I1, P63951, Q3JMZ6, Q7MLS6, Q9JYH8, Q8FNZ0, Q47QS9, Q3BUF7, Q4UU40, Q9RH05, Q3M5D3, Q7P0F6, Q8NNT5, Q31GY9, Q8DTQ7, Q492K3, Q2YZB4, Q57F99, Q81T62, P28736, Q2YPB8, Q5LAZ8, P9WNW9, P59398, P59400, O66976, P60860, P63952, O30027, Q5HSJ3, Q2RQM7, Q9HVW9, P62459, Q8PZR8, P59399, Q87C29, Q9F854, Q87WV5, Q9CC57, Q5P790, Q3AR25, Q8YF59, Q465R4, P59401, Q9X0D1, Q7VSZ1, Q7UX39, Q3SI67, Q63Q86, Q66C49, Q92E84, P57201, Q5H0L1, Q5JFR5, Q8CQ95, Q5NAY4, Q9PM77, Q970Y9, Q4FQF8, Q8TL41, Q6ANR3, Q3KHZ2, Q64RE7, Q6KZD1, Q9RQ82, Q8KEY6, Q67KH6, Q4KI73, Q32EE9
"@
$ws.Range("K2").Value = $k2

# --- Row 3 values ---
$ws.Range("D3").Value = 30.21155927414955
$ws.Range("E3").Value = 0.4688516113516113
$ws.Range("F3").Value = 0.26
$ws.Range("G3").Value = 0.3839285714285714
$ws.Range("H3").Value = 0.7626262626262625
$i3 = @"
Q89QH2, A4YXN2, B3QBS6, Q8XBR7, B7M6P3, A4G242, A5EGD7, Q46S66, Q2IUI7, B1X9P6, Q31Y97, Q6N8F8, B1IX88, A7ZPI2, B1LMH0, Q217M3, Q46S72, A4G241, Q8FFE8, B6JE29, B7NPQ8, C4ZVR1, A7ICK2, B7UG84, Q0TF87, Q82M40, A6T0J2, B2TWX3, Q0K0H8, A9X6P7, O87838, Q3YZF6, B5YYX4, P69902, Q07Q82, Q1R8Z2, B6I6S5, A8A2M8, Q13RQ4, B7LBS7, P69903, O06644, B7MH34, Q32DG9, B7MY33, A1ADQ1, B7N5X4, Q139H7
"@
$ws.Range("I3").Value = $i3
$ws.Range("J3").Value = "O87943, O87942, O87944"
$ws.Range("K3").Value = $k2
